{"js": "// Fill in the first empty log row of the hours-tracking table with the\n// new entry: date, hours and description (login/signup work).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Locate the first completely empty data row (skip the header row).\nlet targetRow = -1;\nfor (let r = 1; r < table.rowCount; r++) {\n  const cell = table.getCell(r, 0);\n  cell.load(\"value\");\n  await context.sync();\n  if (cell.value.trim() === \"\") {\n    targetRow = r;\n    break;\n  }\n}\n\nif (targetRow === -1) {\n  throw new Error(\"No empty row found in the table.\");\n}\n\ntable.getCell(targetRow, 0).value = \"15.2.2023\";\ntable.getCell(targetRow, 1).value = \"3\";\ntable.getCell(targetRow, 2).value =\n  \"Tehtiin yhdess\u00e4 Login/signup ikkunoita. Idea on t\u00e4ss\u00e4 kohtaa selvill\u00e4, mutta ohjelma ei toimi viel\u00e4.\";\n\nawait context.sync();\n", "ps1": "# Fill in the first empty log row of the hours-tracking table with the\n# new entry: date, hours and description (login/signup work).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$targetRow = -1\nfor ($i = 2; $i -le $rowCount; $i++) {\n    $cellText = $t.Cell($i, 1).Range.Text\n    $cellText = $cellText -replace \"[\\x07\\x0d]\", \"\"\n    if ($cellText -eq \"\") {\n        $targetRow = $i\n        break\n    }\n}\n\nif ($targetRow -eq -1) {\n    throw \"No empty row found in the table.\"\n}\n\n$t.Cell($targetRow, 1).Range.Text = \"15.2.2023\"\n$t.Cell($targetRow, 2).Range.Text = \"3\"\n$t.Cell($targetRow, 3).Range.Text = \"Tehtiin yhdess\u00e4 Login/signup ikkunoita. Idea on t\u00e4ss\u00e4 kohtaa selvill\u00e4, mutta ohjelma ei toimi viel\u00e4.\"\n"}
